$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force text format so numeric-looking
# values like "0.9996" / "6.010" are preserved exactly as text,
# matching the source inlineStr cell type instead of being
# auto-converted to a floating point number by Excel.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.188.59'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.917.12'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '329.73'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9996'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5227'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4075'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08522'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '42.88'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.439'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.913.95'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.404'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '95.03'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06686'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.39'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9994'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.010'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '30.208.84'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.211'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.144.43'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '160.83'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.416'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.083'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.016'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.603'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02491'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06578'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.180'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.877'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6541'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.65'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.244'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.27'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.082'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '124.74'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.165'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '79.62'

# Volume(1h) column (E) updates - plain percentage text, safe as-is.
$ws.Range('E2').Value = '  +6.04%  '
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('E4').Value = '  -0.85%  '
$ws.Range('E5').Value = '  +4.60%  '
$ws.Range('E6').Value = '  -0.85%  '
$ws.Range('E7').Value = '  +3.08%  '
$ws.Range('E8').Value = '  +4.25%  '
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('E11').Value = '  +1.88%  '
$ws.Range('E12').Value = '  +10.67%  '
$ws.Range('E13').Value = '  +4.06%  '
$ws.Range('E14').Value = '  +2.83%  '
$ws.Range('E15').Value = '  +2.44%  '
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('E17').Value = '  +4.15%  '
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('E20').Value = '  +4.38%  '
$ws.Range('E21').Value = '  -0.82%  '
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('E23').Value = '  +5.92%  '
$ws.Range('E24').Value = '  +2.38%  '
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('E26').Value = '  +3.19%  '
$ws.Range('E27').Value = '  +1.98%  '
$ws.Range('E28').Value = '  +2.96%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').Value = '  +2.59%  '
$ws.Range('E31').Value = '  +4.68%  '
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('E33').Value = '  +4.81%  '
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('E35').Value = '  +1.86%  '
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('E37').Value = '  +2.32%  '
$ws.Range('E38').Value = '  +4.36%  '
$ws.Range('E39').Value = '  +3.33%  '
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('E41').Value = '  +3.08%  '
$ws.Range('E42').Value = '  +5.24%  '
$ws.Range('E43').Value = '  +0.84%  '
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('E45').Value = '  +1.77%  '
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('E47').Value = '  +4.38%  '
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('E49').Value = '  +2.19%  '
$ws.Range('E50').Value = '  +2.46%  '
$ws.Range('E51').Value = '  +4.62%  '
